$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $c = $ws.Range($cellAddr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "27.408.12"
Set-TextValue "E2" "  +1.25%  "
Set-TextValue "D3" "1.780.37"
Set-TextValue "E3" "  +3.79%  "
Set-TextValue "D4" "1.000"
Set-TextValue "E4" "  -0.06%  "
Set-TextValue "D5" "313.94"
Set-TextValue "E5" "  +1.25%  "
Set-TextValue "E6" "  -0.06%  "
Set-TextValue "D7" "0.5270"
Set-TextValue "E7" "  +9.80%  "
Set-TextValue "D8" "0.3772"
Set-TextValue "E8" "  +9.37%  "
Set-TextValue "D9" "42.84"
Set-TextValue "E9" "  +1.57%  "
Set-TextValue "D10" "0.07414"
Set-TextValue "E10" "  +2.20%  "
Set-TextValue "E11" "  +5.42%  "
Set-TextValue "D12" "1.000"
Set-TextValue "E12" "  -0.06%  "
Set-TextValue "D13" "20.73"
Set-TextValue "E13" "  +5.00%  "
Set-TextValue "D14" "6.112"
Set-TextValue "E14" "  +4.78%  "
Set-TextValue "D15" "1.779.37"
Set-TextValue "E15" "  +3.86%  "
Set-TextValue "D16" "6.986"
Set-TextValue "E16" "  +2.39%  "
Set-TextValue "D17" "89.87"
Set-TextValue "E17" "  +3.19%  "
Set-TextValue "E18" "  +2.23%  "
Set-TextValue "D19" "0.06439"
Set-TextValue "E19" "  +0.94%  "
Set-TextValue "D20" "0.9999"
Set-TextValue "E20" "  -0.05%  "
Set-TextValue "D21" "16.83"
Set-TextValue "E21" "  +2.28%  "
Set-TextValue "D22" "5.901"
Set-TextValue "E22" "  +5.04%  "
Set-TextValue "D23" "27.444.16"
Set-TextValue "E23" "  +1.19%  "
Set-TextValue "E24" "  +4.46%  "
Set-TextValue "E25" "  -0.21%  "
Set-TextValue "D26" "155.55"
Set-TextValue "E26" "  +3.13%  "
Set-TextValue "D27" "20.25"
Set-TextValue "E27" "  +1.16%  "
Set-TextValue "D28" "2.359"
Set-TextValue "E28" "  +14.83%  "
Set-TextValue "D29" "1.986.44"
Set-TextValue "E29" "  +4.05%  "
Set-TextValue "D30" "121.25"
Set-TextValue "E30" "  +0.55%  "
Set-TextValue "E31" "  +5.46%  "
Set-TextValue "D32" "0.1015"
Set-TextValue "E32" "  +9.66%  "
Set-TextValue "D33" "5.607"
Set-TextValue "E33" "  +5.62%  "
Set-TextValue "D34" "3.633"
Set-TextValue "E34" "  +0.98%  "
Set-TextValue "D35" "0.02260"
Set-TextValue "E35" "  +3.80%  "
Set-TextValue "D36" "0.05997"
Set-TextValue "E36" "  +2.78%  "
Set-TextValue "E37" "  +3.92%  "
Set-TextValue "B38" "InternetComputer(DFINITY)"
Set-TextValue "C38" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D38" "4.914"
Set-TextValue "E38" "  +4.37%  "
Set-TextValue "B39" "Algorand"
Set-TextValue "C39" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D39" "0.2051"
Set-TextValue "E39" "  +3.46%  "
Set-TextValue "D40" "0.6138"
Set-TextValue "E40" "  +3.27%  "
Set-TextValue "B41" "FraxShare"
Set-TextValue "C41" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D41" "8.220"
Set-TextValue "E41" "  +9.79%  "
Set-TextValue "B42" "WEMIXTOKEN"
Set-TextValue "C42" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D42" "1.433"
Set-TextValue "E42" "  -2.68%  "
Set-TextValue "D43" "1.134"
Set-TextValue "E43" "  +4.74%  "
Set-TextValue "D44" "13.20"
Set-TextValue "E44" "  +3.56%  "
Set-TextValue "D45" "0.5800"
Set-TextValue "E45" "  +4.34%  "
Set-TextValue "D46" "3.620"
Set-TextValue "E46" "  +1.07%  "
Set-TextValue "D47" "121.84"
Set-TextValue "E47" "  +2.73%  "
Set-TextValue "D48" "1.897"
Set-TextValue "E48" "  +4.18%  "
Set-TextValue "D49" "1.122"
Set-TextValue "E49" "  +2.77%  "
Set-TextValue "D50" "0.06737"
Set-TextValue "E50" "  +1.60%  "
Set-TextValue "D51" "70.99"
Set-TextValue "E51" "  +2.65%  "
